# eims-toi-transect-info.xlsx edits
# - Fix temperature attribute definition text (was incorrectly copy/pasted "salinity")
# - Add Kate Morkeski as a new Personnel entry (metadata Provider) with new NSF award
# - Update NES-LTER Information Manager contact's funding number to the new award
# - Add projectTitle/fundingAgency to Zoe Sandwith's row
# - Restore selections/active cells as recorded for each sheet

$wb = $excel.ActiveWorkbook

# --- Sheet: ColumnHeadersEims (sheet1) ---
$wsEims = $wb.Worksheets.Item("ColumnHeadersEims")
$wsEims.Cells.Item(7,2).Value = "Underway thermosalinograph temperature in degrees Celsius. URI http://vocab.nerc.ac.uk/collection/P01/current/TEMPSZ01/"

# --- Sheet: Personnel (sheet4) ---
$wsPersonnel = $wb.Worksheets.Item("Personnel")

# fill in project/funding columns for Zoe Sandwith row (row 7)
$wsPersonnel.Cells.Item(7,8).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(7,9).Value = "NSF"

# add new personnel row for Kate Morkeski (row 8)
$wsPersonnel.Cells.Item(8,1).Value = "Kate"
$wsPersonnel.Cells.Item(8,3).Value = "Morkeski"
$wsPersonnel.Cells.Item(8,4).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(8,5).Value = "kmorkeski@whoi.edu"
$wsPersonnel.Cells.Item(8,6).Value = "0000-0002-2903-5851"
$wsPersonnel.Cells.Item(8,7).Value = "metadata Provider"
$wsPersonnel.Cells.Item(8,8).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(8,9).Value = "NSF"
$wsPersonnel.Cells.Item(8,10).Value = "OCE-2322676"

$e8 = $wsPersonnel.Range("E8")
$e3 = $wsPersonnel.Range("E3")
$wsPersonnel.Hyperlinks.Add($e8, "mailto:kmorkeski@whoi.edu") | Out-Null
$e8.Style = $e3.Style

# update funding number for NES-LTER Information Manager contact row
$wsPersonnel.Cells.Item(2,10).Value = "OCE-2322676"

# --- selections (must end with Personnel active, matching activeTab=3) ---
$wsEims.Range("B17").Select()

$wsToi = $wb.Worksheets.Item("ColumnHeadersToi")
$wsToi.Range("B10").Select()

$wsPersonnel.Range("J13").Select()
